{"js": "// Append a new 2-row \"bug report\" table right after the existing\n// trailing empty paragraph that follows the last table in the body\n// (i.e. between the two empty paragraphs at the very end of the\n// document, before the final sectPr).\n//\n// The new table mirrors the structure/styling of the preceding bug\n// table (same tcBorders/tcMar/vAlign, bold header row) and its second\n// data cell reproduces a \"pasted from web\" nested-table fragment\n// (two 1x1 helper tables plus a hidden \"vanish\" paragraph) around the\n// \"Default word-splitting behavior for alphanumeric text\" text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The very last body-level paragraph is the final trailing empty\n// paragraph; inserting OOXML \"before\" its range places the new table\n// immediately ahead of it (i.e. right after the other trailing empty\n// paragraph that directly follows the previous table) - exactly where\n// the diff adds it.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst insertionRange = lastParagraph.getRange();\n\nconst tableOoxml =\n  '<w:tbl><w:tblPr><w:tblW w:w=\"0\" w:type=\"auto\"/><w:tblCellMar><w:top w:w=\"15\" w:type=\"dxa\"/><w:left w:w=\"15\" w:type=\"dxa\"/><w:bottom w:w=\"15\" w:type=\"dxa\"/><w:right w:w=\"15\" w:type=\"dxa\"/></w:tblCellMar><w:tblLook w:val=\"04A0\" w:firstRow=\"1\" w:lastRow=\"0\" w:firstColumn=\"1\" w:lastColumn=\"0\" w:noHBand=\"0\" w:noVBand=\"1\"/></w:tblPr><w:tblGrid><w:gridCol w:w=\"622\"/><w:gridCol w:w=\"1695\"/><w:gridCol w:w=\"1166\"/><w:gridCol w:w=\"3857\"/><w:gridCol w:w=\"1676\"/></w:tblGrid><w:tr><w:trPr><w:trHeight w:val=\"640\"/></w:trPr><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Time</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Properties to test</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Approach</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Found failures</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Other notes</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val=\"640\"/></w:trPr><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>12:</w:t></w:r><w:r><w:t>16</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:tbl><w:tblPr><w:tblW w:w=\"0\" w:type=\"auto\"/><w:tblCellSpacing w:w=\"15\" w:type=\"dxa\"/><w:tblCellMar><w:top w:w=\"15\" w:type=\"dxa\"/><w:left w:w=\"15\" w:type=\"dxa\"/><w:bottom w:w=\"15\" w:type=\"dxa\"/><w:right w:w=\"15\" w:type=\"dxa\"/></w:tblCellMar><w:tblLook w:val=\"04A0\" w:firstRow=\"1\" w:lastRow=\"0\" w:firstColumn=\"1\" w:lastColumn=\"0\" w:noHBand=\"0\" w:noVBand=\"1\"/></w:tblPr><w:tblGrid><w:gridCol w:w=\"96\"/></w:tblGrid><w:tr><w:trPr><w:tblCellSpacing w:w=\"15\" w:type=\"dxa\"/></w:trPr><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:vAlign w:val=\"center\"/><w:hideMark/></w:tcPr><w:p/></w:tc></w:tr></w:tbl><w:p><w:pPr><w:rPr><w:vanish/></w:rPr></w:pPr></w:p><w:tbl><w:tblPr><w:tblW w:w=\"0\" w:type=\"auto\"/><w:tblCellSpacing w:w=\"15\" w:type=\"dxa\"/><w:tblCellMar><w:top w:w=\"15\" w:type=\"dxa\"/><w:left w:w=\"15\" w:type=\"dxa\"/><w:bottom w:w=\"15\" w:type=\"dxa\"/><w:right w:w=\"15\" w:type=\"dxa\"/></w:tblCellMar><w:tblLook w:val=\"04A0\" w:firstRow=\"1\" w:lastRow=\"0\" w:firstColumn=\"1\" w:lastColumn=\"0\" w:noHBand=\"0\" w:noVBand=\"1\"/></w:tblPr><w:tblGrid><w:gridCol w:w=\"1655\"/></w:tblGrid><w:tr><w:trPr><w:tblCellSpacing w:w=\"15\" w:type=\"dxa\"/></w:trPr><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:vAlign w:val=\"center\"/><w:hideMark/></w:tcPr><w:p><w:r><w:t xml:space=\"preserve\">Default word-splitting </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>behavior</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> for alphanumeric text</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>Black-box unit testing using Jest.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>words(\\'item1 item2a 123abc\\') splits incorrectly: returns [\"item\",\"1\",\"item\",\"2\",\"a\",\"123\",\"abc\"] instead of [\"item1\",\"item2a\",\"123abc\"]</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>Suspected cause:</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>ASCII regex (</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>reAsciiWord</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>) treats digits and letters as separate segments, not merged tokens.</w:t></w:r></w:p></w:tc></w:tr></w:tbl>';\n\n// insertOoxml on a Word.js object needs the full OOXML \"flat\" package\n// wrapper (the same shape Range.WordOpenXML / the clipboard use), not\n// a bare WordprocessingML fragment.\nconst packageOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + tableOoxml + '<w:sectPr/></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ninsertionRange.insertOoxml(packageOoxml, Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document ends with two empty paragraphs (the one right after the\n# last existing table, and the final one before sectPr). The new table\n# must be inserted between them, leaving both paragraphs themselves\n# untouched.\n$lastPara = $d.Paragraphs.Last\n$firstTrailingPara = $lastPara.Previous()\n\n# Create a fresh empty paragraph right after the first trailing empty\n# paragraph; this reserves the exact insertion slot without disturbing\n# either of the original empty paragraphs.\n$insertionPoint = $firstTrailingPara.Range.Duplicate\n$insertionPoint.Collapse(0)\n$insertionPoint.InsertParagraphAfter()\n$placeholderPara = $firstTrailingPara.Next()\n\n$tableXml = @'\n<w:tbl><w:tblPr><w:tblW w:w=\"0\" w:type=\"auto\"/><w:tblCellMar><w:top w:w=\"15\" w:type=\"dxa\"/><w:left w:w=\"15\" w:type=\"dxa\"/><w:bottom w:w=\"15\" w:type=\"dxa\"/><w:right w:w=\"15\" w:type=\"dxa\"/></w:tblCellMar><w:tblLook w:val=\"04A0\" w:firstRow=\"1\" w:lastRow=\"0\" w:firstColumn=\"1\" w:lastColumn=\"0\" w:noHBand=\"0\" w:noVBand=\"1\"/></w:tblPr><w:tblGrid><w:gridCol w:w=\"622\"/><w:gridCol w:w=\"1695\"/><w:gridCol w:w=\"1166\"/><w:gridCol w:w=\"3857\"/><w:gridCol w:w=\"1676\"/></w:tblGrid><w:tr><w:trPr><w:trHeight w:val=\"640\"/></w:trPr><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Time</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Properties to test</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Approach</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Found failures</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Other notes</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val=\"640\"/></w:trPr><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>12:</w:t></w:r><w:r><w:t>16</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:tbl><w:tblPr><w:tblW w:w=\"0\" w:type=\"auto\"/><w:tblCellSpacing w:w=\"15\" w:type=\"dxa\"/><w:tblCellMar><w:top w:w=\"15\" w:type=\"dxa\"/><w:left w:w=\"15\" w:type=\"dxa\"/><w:bottom w:w=\"15\" w:type=\"dxa\"/><w:right w:w=\"15\" w:type=\"dxa\"/></w:tblCellMar><w:tblLook w:val=\"04A0\" w:firstRow=\"1\" w:lastRow=\"0\" w:firstColumn=\"1\" w:lastColumn=\"0\" w:noHBand=\"0\" w:noVBand=\"1\"/></w:tblPr><w:tblGrid><w:gridCol w:w=\"96\"/></w:tblGrid><w:tr><w:trPr><w:tblCellSpacing w:w=\"15\" w:type=\"dxa\"/></w:trPr><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:vAlign w:val=\"center\"/><w:hideMark/></w:tcPr><w:p/></w:tc></w:tr></w:tbl><w:p><w:pPr><w:rPr><w:vanish/></w:rPr></w:pPr></w:p><w:tbl><w:tblPr><w:tblW w:w=\"0\" w:type=\"auto\"/><w:tblCellSpacing w:w=\"15\" w:type=\"dxa\"/><w:tblCellMar><w:top w:w=\"15\" w:type=\"dxa\"/><w:left w:w=\"15\" w:type=\"dxa\"/><w:bottom w:w=\"15\" w:type=\"dxa\"/><w:right w:w=\"15\" w:type=\"dxa\"/></w:tblCellMar><w:tblLook w:val=\"04A0\" w:firstRow=\"1\" w:lastRow=\"0\" w:firstColumn=\"1\" w:lastColumn=\"0\" w:noHBand=\"0\" w:noVBand=\"1\"/></w:tblPr><w:tblGrid><w:gridCol w:w=\"1655\"/></w:tblGrid><w:tr><w:trPr><w:tblCellSpacing w:w=\"15\" w:type=\"dxa\"/></w:trPr><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:vAlign w:val=\"center\"/><w:hideMark/></w:tcPr><w:p><w:r><w:t xml:space=\"preserve\">Default word-splitting </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>behavior</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> for alphanumeric text</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>Black-box unit testing using Jest.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>words('item1 item2a 123abc') splits incorrectly: returns [\"item\",\"1\",\"item\",\"2\",\"a\",\"123\",\"abc\"] instead of [\"item1\",\"item2a\",\"123abc\"]</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w=\"0\" w:type=\"auto\"/><w:tcBorders><w:top w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:left w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:bottom w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/><w:right w:val=\"single\" w:sz=\"4\" w:space=\"0\" w:color=\"000000\"/></w:tcBorders><w:tcMar><w:top w:w=\"20\" w:type=\"dxa\"/><w:left w:w=\"20\" w:type=\"dxa\"/><w:bottom w:w=\"100\" w:type=\"dxa\"/><w:right w:w=\"20\" w:type=\"dxa\"/></w:tcMar><w:vAlign w:val=\"bottom\"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>Suspected cause:</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>ASCII regex (</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>reAsciiWord</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>) treats digits and letters as separate segments, not merged tokens.</w:t></w:r></w:p></w:tc></w:tr></w:tbl>\n'@\n\n# Replace the placeholder paragraph's range with the new table's OOXML.\n$placeholderPara.Range.InsertXML($tableXml)\n"}
